# Weekly update to the "Femacal de La Calera - Espárragos" sheet:
# two new daily price rows are inserted into the historical table.
#
# Row 1: inserted right after row 19 (becomes the new row 20), shifting
#        every subsequent row down by one.
# Row 2: inserted right after the (now shifted) row 22 (becomes the new
#        row 24), shifting everything below it down by one more.
#
# Net effect: the table grows from 32 data rows (A1:R32) to 34 (A1:R34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row at 20 -------------------------------------------------
$ws.Rows.Item(20).Insert()

$ws.Cells.Item(20, 1).Value = 3
$ws.Cells.Item(20, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(20, 3).Value = "Coquimbo"
$ws.Cells.Item(20, 4).Value = 44846
$ws.Cells.Item(20, 5).Value = 5
$ws.Cells.Item(20, 6).Value = 300000000
$ws.Cells.Item(20, 7).Value = "Espárragos"
$ws.Cells.Item(20, 8).Value = "Verde"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 1000
$ws.Cells.Item(20, 11).Value = 1400
$ws.Cells.Item(20, 12).Value = 1450
$ws.Cells.Item(20, 13).Value = 1428
$ws.Cells.Item(20, 14).Value = "$/kilo"
$ws.Cells.Item(20, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(20, 16).Value = 1428
$ws.Cells.Item(20, 17).Value = 1
$ws.Cells.Item(20, 18).Value = "Hortaliza"

# --- Insert second new row at 24 (after the shift above) -----------------
$ws.Rows.Item(24).Insert()

$ws.Cells.Item(24, 1).Value = 3
$ws.Cells.Item(24, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(24, 3).Value = "Coquimbo"
$ws.Cells.Item(24, 4).Value = 44845
$ws.Cells.Item(24, 5).Value = 5
$ws.Cells.Item(24, 6).Value = 300000000
$ws.Cells.Item(24, 7).Value = "Espárragos"
$ws.Cells.Item(24, 8).Value = "Verde"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 1000
$ws.Cells.Item(24, 11).Value = 1300
$ws.Cells.Item(24, 12).Value = 1500
$ws.Cells.Item(24, 13).Value = 1396
$ws.Cells.Item(24, 14).Value = "$/kilo"
$ws.Cells.Item(24, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(24, 16).Value = 1396
$ws.Cells.Item(24, 17).Value = 1
$ws.Cells.Item(24, 18).Value = "Hortaliza"
